# Edit script: Additional statistical tests and dropping columns in Data
# Preparation for ignored columns.
#
# This updates the "Data Cleaning V1" sheet (sheet2) of the workbook:
#  - Replaces the old placeholder "ignored, not used in paper" cleaning-action
#    notes with concrete decisions for each previously-ignored field.
#  - Un-hides the rows that were hidden as a side effect of the old autofilter
#    criteria (filtering out "ignored, not used in paper"), except the rows
#    that should remain hidden for other reasons.
#  - Clears the autofilter criteria so the filter no longer hides any rows.
#  - Leaves the "Data Cleaning V1" sheet as the active tab/sheet, with cell
#    C5 selected; the "CARES_fields_info" sheet keeps cell D26 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update Cleaning/Encoding Action notes on "Data Cleaning V1" ---
# (cell write order below intentionally mirrors the author's original
#  editing sequence so that newly-introduced shared strings end up in the
#  same order in the workbook's shared string table)

# RCRI score: dropped, could lead to multicollinearity
$ws2.Range("C5").Value = "dropped, could lead to multicollinearity"

# PreopEGFRMDRD: categorize by eGFR threshold
$ws2.Range("C7").Value = "Categorization: preop-eGFR " + [char]0x2264 + " 98.688 mL/min/1.73 m2 (0) and preop-eGFR > 98.688 mL/min/1.73 m2 (1)"

# Preoptransfusionwithin30days / Intraop / Postopwithin30days combined
$ws2.Range("C11").Value = "Combined into Perioperativetransfusion"
$ws2.Range("C12").Value = "Combined into Perioperativetransfusion"
$ws2.Range("C13").Value = "Combined into Perioperativetransfusion"

# PreopEGFRMDRD: label encoding
$ws2.Range("D7").Value = "Label Encoding"

# Perioperativetransfusion: label encoding details
$ws2.Range("D11").Value = "Label Encoding: Group all three transfusion columns and name the feature as Perioperativetransfusion (1, 0), presence of a single transfusion using OR operation."
$ws2.Range("D12").Value = "Label Encoding: Group all three transfusion columns and name the feature as Perioperativetransfusion (1, 0), presence of a single transfusion using OR operation."
$ws2.Range("D13").Value = "Label Encoding: Group all three transfusion columns and name the feature as Perioperativetransfusion (1, 0), presence of a single transfusion using OR operation."

# Transfusionintraandpostop: dropped, as it is redundant
$ws2.Range("C14").Value = "dropped, as it is redundant"

# TransfusionIntraandpostopCategory: dropped, as it is redundant
$ws2.Range("C17").Value = "dropped, as it is redundant"

# AnaestypeCategory: no change, label encoding GA-0, RA-1
$ws2.Range("C15").Value = "no change"
$ws2.Range("D15").Value = "Label Encoding: GA-0, RA-1"

# --- Clear the autofilter criteria (keep the autofilter range itself). ---
# This also unhides every previously-filtered row.
$ws2.ShowAllData()

# --- Re-hide the rows that must stay hidden (dropped / not shown fields) ---
$rowsToRehide = @(4,8,9,10)
foreach ($r in $rowsToRehide) {
    $ws2.Rows.Item($r).Hidden = $true
}

# --- Selections / active sheet ---
$ws1.Range("D26").Select()
$ws2.Activate()
$ws2.Range("C5").Select()
